$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 29 ("Provincia de Quillota",
# 2022-10-21), pushing the existing rows 29-36 down to rows 30-37.
$ws.Rows("29:29").Insert()

$ws.Range("A29").Value = 3
$ws.Range("B29").Value = "Femacal de La Calera"
$ws.Range("C29").Value = "Coquimbo"
$ws.Range("D29").Value = 44855
$ws.Range("E29").Value = 5
$ws.Range("F29").Value = 300000000
$ws.Range("G29").Value = "Espárragos"
$ws.Range("H29").Value = "Verde"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 120
$ws.Range("K29").Value = 1400
$ws.Range("L29").Value = 1400
$ws.Range("M29").Value = 1400
$ws.Range("N29").Value = "$/kilo"
$ws.Range("O29").Value = "Provincia de Quillota"
$ws.Range("P29").Value = 1400
$ws.Range("Q29").Value = 1
$ws.Range("R29").Value = "Hortaliza"
